$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right column 4 -> 5, Wrong column -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total): Right column 88 -> 110, Wrong column 0 -> -0
$ws.Range("B12").Value = 110
$ws.Range("C12").Value = -0

# E12 label update to match new totals
$ws.Range("E12").Value = "110.0/140"
